# Apply updated cryptocurrency price/volume data (and two row swaps)
# to match the latest scrape, as produced by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'29.500.51"
$ws.Range("E2").Formula = "'  -3.10%  "
$ws.Range("D3").Formula = "'1.994.64"
$ws.Range("E3").Formula = "'  -6.31%  "
$ws.Range("D4").Formula = "'1.008"
$ws.Range("E4").Formula = "'  +0.29%  "
$ws.Range("D5").Formula = "'329.69"
$ws.Range("E5").Formula = "'  -5.32%  "
$ws.Range("D7").Formula = "'0.5000"
$ws.Range("E7").Formula = "'  -4.35%  "
$ws.Range("D8").Formula = "'0.4218"
$ws.Range("E8").Formula = "'  -6.15%  "
$ws.Range("D9").Formula = "'51.89"
$ws.Range("E9").Formula = "'  -4.18%  "
$ws.Range("D10").Formula = "'0.08895"
$ws.Range("E10").Formula = "'  -5.37%  "
$ws.Range("D11").Formula = "'1.121"
$ws.Range("E11").Formula = "'  -5.66%  "
$ws.Range("D12").Formula = "'23.34"
$ws.Range("E12").Formula = "'  -8.59%  "
$ws.Range("D13").Formula = "'8.118"
$ws.Range("E13").Formula = "'  -6.91%  "
$ws.Range("D14").Formula = "'1.984.45"
$ws.Range("E14").Formula = "'  -6.94%  "
$ws.Range("D15").Formula = "'6.513"
$ws.Range("E15").Formula = "'  -6.75%  "
$ws.Range("D16").Formula = "'96.23"
$ws.Range("E16").Formula = "'  -6.91%  "
$ws.Range("D17").Formula = "'1.008"
$ws.Range("E17").Formula = "'  +0.22%  "
$ws.Range("E18").Formula = "'  -5.72%  "
$ws.Range("D19").Formula = "'0.06636"
$ws.Range("E19").Formula = "'  -1.20%  "
$ws.Range("D20").Formula = "'19.74"
$ws.Range("E20").Formula = "'  -8.75%  "
$ws.Range("E21").Formula = "'  +0.21%  "
$ws.Range("D22").Formula = "'5.971"
$ws.Range("E22").Formula = "'  -5.90%  "
$ws.Range("D23").Formula = "'29.526.44"
$ws.Range("E23").Formula = "'  -3.06%  "
$ws.Range("D24").Formula = "'11.89"
$ws.Range("E24").Formula = "'  -7.07%  "
$ws.Range("D25").Formula = "'2.276"
$ws.Range("E25").Formula = "'  -2.67%  "
$ws.Range("D26").Formula = "'157.97"
$ws.Range("E26").Formula = "'  -3.53%  "
$ws.Range("B27").Formula = "'InternetComputer(DFINITY)"
$ws.Range("C27").Formula = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Formula = "'6.571"
$ws.Range("E27").Formula = "'  -4.95%  "
$ws.Range("B28").Formula = "'EthereumClassic"
$ws.Range("C28").Formula = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Formula = "'20.61"
$ws.Range("E28").Formula = "'  -7.71%  "
$ws.Range("D29").Formula = "'2.335"
$ws.Range("E29").Formula = "'  -8.84%  "
$ws.Range("E30").Formula = "'  -5.00%  "
$ws.Range("D31").Formula = "'1.053"
$ws.Range("E31").Formula = "'  -10.07%  "
$ws.Range("D32").Formula = "'0.09946"
$ws.Range("E32").Formula = "'  -6.32%  "
$ws.Range("D33").Formula = "'1.559"
$ws.Range("E33").Formula = "'  -13.97%  "
$ws.Range("D34").Formula = "'5.840"
$ws.Range("E34").Formula = "'  -7.65%  "
$ws.Range("D35").Formula = "'3.788"
$ws.Range("E35").Formula = "'  -4.26%  "
$ws.Range("D36").Formula = "'9.598"
$ws.Range("E36").Formula = "'  -10.76%  "
$ws.Range("D37").Formula = "'0.02462"
$ws.Range("E37").Formula = "'  -7.13%  "
$ws.Range("D38").Formula = "'0.06359"
$ws.Range("E38").Formula = "'  -7.64%  "
$ws.Range("D39").Formula = "'1.289"
$ws.Range("E39").Formula = "'  -3.78%  "
$ws.Range("D40").Formula = "'0.6533"
$ws.Range("E40").Formula = "'  -8.88%  "
$ws.Range("D41").Formula = "'11.73"
$ws.Range("E41").Formula = "'  -8.52%  "
$ws.Range("D42").Formula = "'0.2067"
$ws.Range("E42").Formula = "'  -8.61%  "
$ws.Range("D43").Formula = "'1.007"
$ws.Range("E43").Formula = "'  +0.27%  "
$ws.Range("D44").Formula = "'0.6352"
$ws.Range("E44").Formula = "'  -8.96%  "
$ws.Range("B45").Formula = "'EnergySwap"
$ws.Range("C45").Formula = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Formula = "'13.47"
$ws.Range("E45").Formula = "'  -8.79%  "
$ws.Range("B46").Formula = "'NEARProtocol"
$ws.Range("C46").Formula = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Formula = "'2.215"
$ws.Range("E46").Formula = "'  -7.82%  "
$ws.Range("D47").Formula = "'1.272"
$ws.Range("E47").Formula = "'  -0.47%  "
$ws.Range("D48").Formula = "'3.525"
$ws.Range("E48").Formula = "'  -3.21%  "
$ws.Range("D49").Formula = "'0.00000000337"
$ws.Range("E49").Formula = "'  -2.92%  "
$ws.Range("D50").Formula = "'0.07013"
$ws.Range("E50").Formula = "'  -3.05%  "
$ws.Range("E51").Formula = "'  -6.53%  "
